$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-interpreted by Excel as a
# number (e.g. "252.39") need to be explicitly formatted as Text first so
# the literal string is preserved, matching the source data which stores
# these as plain text.
$ws.Range('D2').Value = '37.566.71'
$ws.Range('E2').Value = '  +6.26%  '
$ws.Range('D3').Value = '2.054.55'
$ws.Range('E3').Value = '  +3.42%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '252.39'
$ws.Range('E5').Value = '  +4.92%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.648'
$ws.Range('E6').Value = '  +2.67%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '66.57'
$ws.Range('E7').Value = '  +18.69%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  +6.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '59.44'
$ws.Range('E10').Value = '  +0.36%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0755'
$ws.Range('E11').Value = '  +4.53%  '
$ws.Range('E12').Value = '  +1.16%  '
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.914'
$ws.Range('E13').Value = '  +2.30%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.28'
$ws.Range('E14').Value = '  +7.16%  '
$ws.Range('E15').Value = '  +3.43%  '
$ws.Range('E16').Value = '  +7.40%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '20.72'
$ws.Range('E17').Value = '  +21.60%  '
$ws.Range('D18').Value = '2.053.52'
$ws.Range('E18').Value = '  +3.39%  '
$ws.Range('D19').Value = '37.536.03'
$ws.Range('E19').Value = '  +6.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '73.46'
$ws.Range('E20').Value = '  +5.22%  '
$ws.Range('D21').Value = '0.0₃0879'
$ws.Range('E21').Value = '  +5.70%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.36'
$ws.Range('E22').Value = '  +7.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.83'
$ws.Range('E23').Value = '  +2.81%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.74'
$ws.Range('E24').Value = '  +21.81%  '
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('E26').Value = '  +4.65%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.59'
$ws.Range('E27').Value = '  +5.50%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '165.58'
$ws.Range('E28').Value = '  +1.84%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.92'
$ws.Range('E29').Value = '  +2.86%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.24'
$ws.Range('E30').Value = '  +10.55%  '
$ws.Range('E31').Value = '  +3.14%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.22'
$ws.Range('E32').Value = '  +7.81%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.111'
$ws.Range('E33').Value = '  +24.21%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.74'
$ws.Range('E34').Value = '  +11.97%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0616'
$ws.Range('E35').Value = '  +5.78%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.45'
$ws.Range('E36').Value = '  +10.59%  '
$ws.Range('E37').Value = '  +0.15%  '
$ws.Range('E38').Value = '  +25.05%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.82'
$ws.Range('E39').Value = '  +1.09%  '
$ws.Range('E40').Value = '  +17.48%  '
$ws.Range('E41').Value = '  +5.42%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.94'
$ws.Range('E42').Value = '  +4.82%  '
$ws.Range('E43').Value = '  +6.06%  '
$ws.Range('E44').Value = '  +6.49%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.14'
$ws.Range('E45').Value = '  +10.07%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '17.08'
$ws.Range('E46').Value = '  +11.27%  '
$ws.Range('E47').Value = '  +20.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '95.53'
$ws.Range('E48').Value = '  +6.62%  '
$ws.Range('D49').Value = '1.430.03'
$ws.Range('E49').Value = '  +4.82%  '
$ws.Range('E50').Value = '  +2.50%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '47.63'
$ws.Range('E51').Value = '  +5.63%  '
